$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 5
$ws.Range("H5").Value = 142857620
$ws.Range("J5").Value = 250000700
$ws.Range("L5").Value = 250000700
$ws.Range("N5").Value = -250000930

# Row 9
$ws.Range("H9").Value = 157425.86
$ws.Range("J9").Value = 1175
$ws.Range("L9").Value = 1175
$ws.Range("N9").Value = -1513

# Row 53
$ws.Range("H53").Value = 1152.4706
$ws.Range("I53").Value = 2538.75
$ws.Range("J53").Value = 725.9231
$ws.Range("K53").Value = 2538.75
$ws.Range("L53").Value = 725.9231
$ws.Range("M53").Value = -1901.75
$ws.Range("N53").Value = -1999.9231

# Row 80
$ws.Range("H80").Value = 332.25
$ws.Range("I80").Value = 140.33333
$ws.Range("J80").Value = 620.125
$ws.Range("K80").Value = 420.99999
$ws.Range("L80").Value = 1860.375
$ws.Range("M80").Value = 577.00001
$ws.Range("N80").Value = -3856.375

# Row 83
$ws.Range("H83").Value = 332.25
$ws.Range("I83").Value = 140.33333
$ws.Range("J83").Value = 620.125
$ws.Range("K83").Value = 1262.99997
$ws.Range("L83").Value = 5581.125
$ws.Range("M83").Value = 3729.00003
$ws.Range("N83").Value = -15565.125

# Row 92
$ws.Range("H92").Value = 270.78946
$ws.Range("I92").Value = 245.84616
$ws.Range("K92").Value = 245.84616
$ws.Range("M92").Value = 1002.15384

# Row 100
$ws.Range("H100").Value = 11935.5
$ws.Range("I100").Value = 329.5
$ws.Range("J100").Value = 17738.5
$ws.Range("K100").Value = 329.5
$ws.Range("L100").Value = 17738.5
$ws.Range("M100").Value = 211.5
$ws.Range("N100").Value = -18820.5

# Row 116
$ws.Range("H116").Value = 9973.429
$ws.Range("I116").Value = 9973.429
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 9973.429
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -6531.429
$ws.Range("N116").ClearContents()

# Row 125
$ws.Range("H125").Value = 2709.8462
$ws.Range("I125").Value = 1438.7142
$ws.Range("K125").Value = 12948.4278
$ws.Range("M125").Value = -10488.4278

# Row 137
$ws.Range("H137").Value = 996.9394
$ws.Range("I137").Value = 793.0714
$ws.Range("K137").Value = 2379.2142
$ws.Range("M137").Value = 170.7857999999997

# Row 138
$ws.Range("H138").Value = 2879.3953
$ws.Range("I138").Value = 4416.864
$ws.Range("J138").Value = 2350.8906
$ws.Range("K138").Value = 13250.592
$ws.Range("L138").Value = 7052.6718
$ws.Range("M138").Value = -8110.591999999999
$ws.Range("N138").Value = -17332.6718

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 773.13794
$ws.Range("J2").Value = 1649.75
$ws.Range("L2").Value = 1649.75
$ws.Range("N2").Value = -1875.75

# Row 32
$ws.Range("H32").Value = 6284.976
$ws.Range("I32").Value = 5938.268
$ws.Range("K32").Value = 5938.268
$ws.Range("M32").Value = -5651.268

# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# Row 45
$ws.Range("H45").Value = 3035.9644
$ws.Range("I45").Value = 1886.04
$ws.Range("J45").Value = 3963.3225
$ws.Range("K45").Value = 1886.04
$ws.Range("L45").Value = 3963.3225
$ws.Range("M45").Value = -1509.04
$ws.Range("N45").Value = -4717.3225

# Row 97
$ws.Range("H97").Value = 2252.149
$ws.Range("I97").Value = 489.94446
$ws.Range("J97").Value = 8019.364
$ws.Range("K97").Value = 489.94446
$ws.Range("L97").Value = 8019.364
$ws.Range("M97").Value = 6.055540000000008
$ws.Range("N97").Value = -9011.364

# Row 116
$ws.Range("H116").Value = 773.13794
$ws.Range("J116").Value = 1649.75
$ws.Range("L116").Value = 1649.75
$ws.Range("N116").Value = -6237.75

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 773.13794
$ws.Range("J3").Value = 1649.75
$ws.Range("L3").Value = 1649.75
$ws.Range("N3").Value = -1877.75

# Row 134
$ws.Range("H134").Value = 8772.23
$ws.Range("I134").Value = 9115.655000000001
$ws.Range("J134").Value = 7776.3
$ws.Range("K134").Value = 27346.965
$ws.Range("L134").Value = 23328.9
$ws.Range("M134").Value = -24811.965
$ws.Range("N134").Value = -28398.9

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 113
$ws.Range("H113").Value = 1290
$ws.Range("J113").Value = 1276.1111
$ws.Range("L113").Value = 3828.3333
$ws.Range("N113").Value = -8168.3333

# Row 126
$ws.Range("H126").Value = 20697.8
$ws.Range("I126").Value = 11395.6
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 34186.8
$ws.Range("L126").Value = 90000
$ws.Range("M126").Value = -29246.8
$ws.Range("N126").Value = -99880

# Row 131
$ws.Range("H131").Value = 12381043
$ws.Range("I131").Value = 7408198.5
$ws.Range("K131").Value = 22224595.5
$ws.Range("M131").Value = -22219555.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 82
$ws.Range("H82").Value = 1009.9091
$ws.Range("I82").Value = 892.2
$ws.Range("J82").Value = 1108
$ws.Range("K82").Value = 892.2
$ws.Range("L82").Value = 1108
$ws.Range("M82").Value = -531.2
$ws.Range("N82").Value = -1830

# Row 85
$ws.Range("H85").Value = 1009.9091
$ws.Range("I85").Value = 892.2
$ws.Range("J85").Value = 1108
$ws.Range("K85").Value = 892.2
$ws.Range("L85").Value = 1108
$ws.Range("M85").Value = 355.8
$ws.Range("N85").Value = -3604

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

# Row 122
$ws.Range("H122").Value = 3383.7693
$ws.Range("J122").Value = 3009.524
$ws.Range("L122").Value = 9028.572
$ws.Range("N122").Value = -13928.572

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 49
$ws.Range("H49").Value = 2300
$ws.Range("I49").Value = 2300
$ws.Range("K49").Value = 2300
$ws.Range("M49").Value = -2070

# Row 54
$ws.Range("H54").Value = 36210.777
$ws.Range("I54").Value = 34724.5
$ws.Range("J54").Value = 37399.8
$ws.Range("K54").Value = 34724.5
$ws.Range("L54").Value = 37399.8
$ws.Range("M54").Value = -34204.5
$ws.Range("N54").Value = -38439.8

# Row 81
$ws.Range("H81").Value = 8703.166999999999
$ws.Range("I81").Value = 10404
$ws.Range("J81").Value = 3600.6667
$ws.Range("K81").Value = 20808
$ws.Range("L81").Value = 7201.3334
$ws.Range("M81").Value = -19747
$ws.Range("N81").Value = -9323.3334

# Row 84
$ws.Range("H84").Value = 8703.166999999999
$ws.Range("I84").Value = 10404
$ws.Range("J84").Value = 3600.6667
$ws.Range("K84").Value = 104040
$ws.Range("L84").Value = 36006.667
$ws.Range("M84").Value = -98736
$ws.Range("N84").Value = -46614.667

# Row 136
$ws.Range("H136").Value = 6646.5405
$ws.Range("I136").Value = 4994.1665
$ws.Range("K136").Value = 14982.4995
$ws.Range("M136").Value = -12432.4995

# Row 137
$ws.Range("H137").Value = 54815.715
$ws.Range("J137").Value = 54815.715
$ws.Range("L137").Value = 54815.715
$ws.Range("N137").Value = -65015.715

Write-Output "Applied Cerberus_Profits market-price update across ALC, ARM, BSM, CUL, LTW, WVR sheets."